$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest quarter column (D: "6 ماهه منتهی به 1399/06"); this shifts
# all subsequent quarter columns one position to the left (E->D, F->E, ... M->L).
$ws.Columns("D:D").Delete()

# Column M is now empty; populate it with the new latest-quarter data
# ("12 ماهه منتهی به 1401/12", published 1402-02-25 (8)).

# Match the column width of the previous column M (now L) for the new M column.
$ws.Range("M1").ColumnWidth = $ws.Range("L1").ColumnWidth

# Row 8: period header
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publication date label (keeps the parenthetical revision counter)
$ws.Range("M9").Value = "1402-02-25 (8)"

# Row 9 column I gets the revised publish-date label too (was shifted from J)
$ws.Range("I9").Value = "1402-02-25 (8)"

# The standalone "1402-02-25" date label (no parenthetical) must be forced to
# text, otherwise Excel parses it as a date serial number.
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "1402-02-25"

# Row 11: Sales (فروش)
$ws.Range("M11").Value = 27260727

# Row 12: Cost of goods sold
$ws.Range("M12").Value = -15164059

# Row 13: Gross profit
$ws.Range("M13").Value = 12096668

# Row 14: G&A expenses
$ws.Range("M14").Value = -4765322

# Row 15 stays 0
$ws.Range("M15").Value = 0

# Row 16 stays 0
$ws.Range("M16").Value = 0

# Row 17: Operating profit
$ws.Range("M17").Value = 7331346

# Row 18: Financial expenses
$ws.Range("M18").Value = -1038607

# Row 19: Net other non-operating income (expenses)
$ws.Range("M19").Value = 101323

# Row 20: Net profit before tax from continuing operations
$ws.Range("M20").Value = 6394062

# Row 21: Tax
$ws.Range("M21").Value = -256002

# Row 22: Net profit from continuing operations
$ws.Range("M22").Value = 6138060

# Row 23 stays 0
$ws.Range("M23").Value = 0

# Row 24: Net profit
$ws.Range("M24").Value = 6138060

# Row 25: EPS after tax
$ws.Range("M25").Value = 1535

# Row 26: Capital
$ws.Range("M26").Value = 4000000

# Row 27: EPS based on latest capital
$ws.Range("M27").Value = 1535

$wb.Save()
